# Update the "classFields" sheet so that field name / field type pairs are
# reshuffled to reflect the regenerated structure report (standard
# relationship between microservices / MSM measure work).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# com.macro.mall.auth.constant.MessageConstant field name relabeling
$ws.Cells.Item(2, 2).Value = "ACCOUNT_DISABLED"
$ws.Cells.Item(3, 2).Value = "ACCOUNT_EXPIRED"
$ws.Cells.Item(4, 2).Value = "USERNAME_PASSWORD_ERROR"
$ws.Cells.Item(5, 2).Value = "ACCOUNT_LOCKED"
$ws.Cells.Item(7, 2).Value = "LOGIN_SUCCESS"

# com.macro.mall.auth.config.Oauth2ServerConfig fields
$ws.Cells.Item(10, 2).Value = "passwordEncoder"
$ws.Cells.Item(10, 4).Value = "org.springframework.security.crypto.password.PasswordEncoder"
$ws.Cells.Item(11, 2).Value = "userDetailsService"
$ws.Cells.Item(11, 4).Value = "com.macro.mall.auth.service.impl.UserServiceImpl"
$ws.Cells.Item(12, 2).Value = "jwtTokenEnhancer"
$ws.Cells.Item(12, 4).Value = "com.macro.mall.auth.component.JwtTokenEnhancer"
$ws.Cells.Item(13, 2).Value = "authenticationManager"
$ws.Cells.Item(13, 4).Value = "org.springframework.security.authentication.AuthenticationManager"

# com.macro.mall.auth.domain.Oauth2TokenDto fields
$ws.Cells.Item(15, 2).Value = "token"
$ws.Cells.Item(16, 2).Value = "refreshToken"
$ws.Cells.Item(17, 2).Value = "tokenHead"

# com.macro.mall.auth.service.impl.UserServiceImpl fields
$ws.Cells.Item(18, 2).Value = "memberService"
$ws.Cells.Item(18, 4).Value = "com.macro.mall.auth.service.UmsMemberService"
$ws.Cells.Item(19, 2).Value = "adminService"
$ws.Cells.Item(19, 4).Value = "com.macro.mall.auth.service.UmsAdminService"
$ws.Cells.Item(20, 2).Value = "request"
$ws.Cells.Item(20, 4).Value = "javax.servlet.http.HttpServletRequest"

# com.macro.mall.auth.domain.Oauth2TokenDto$Oauth2TokenDtoBuilder fields
$ws.Cells.Item(21, 2).Value = "token"
$ws.Cells.Item(22, 2).Value = "refreshToken"
$ws.Cells.Item(23, 2).Value = "tokenHead"
$ws.Cells.Item(23, 4).Value = "java.lang.String"
$ws.Cells.Item(24, 2).Value = "expiresIn"
$ws.Cells.Item(24, 4).Value = "int"

# com.macro.mall.auth.domain.SecurityUser fields
$ws.Cells.Item(25, 2).Value = "password"
$ws.Cells.Item(26, 2).Value = "clientId"
$ws.Cells.Item(26, 4).Value = "java.lang.String"
$ws.Cells.Item(28, 2).Value = "id"
$ws.Cells.Item(28, 4).Value = "java.lang.Long"
